# Generate Report for Handoff
#
# The localization report was re-generated: the "zh-cn" and "de-de"
# localization jobs moved from "Ready for handoff" into "In Translation",
# and the handoff/generation timestamps advanced a few seconds. Excel
# re-fit the "Status" columns to the new (shorter) text, so their
# stored widths shrink too.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("G2").Value = "2016-10-17 11:50:46"

$wsOverview.Columns.Item(5).ColumnWidth = 13.4101848602295
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101848602295

# --- zh-cn sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("H2").Value = "2016-10-17 11:50:36"

$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101848602295

# --- de-de sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("H2").Value = "2016-10-17 11:50:46"

$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101848602295
